$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.754623889923096
$ws.Range("B1").Value = 2.406201601028442
$ws.Range("C1").Value = 2.621898174285889
$ws.Range("D1").Value = 3.381045818328857
$ws.Range("E1").Value = 1.316743612289429
